$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Add Hs (column E) values for Singapore archetypes, row by row (row 1 is the header)
$hsValues = @{
    2  = 0.25   # MULTI_RES
    3  = 0.25   # SINGLE_RES
    4  = 0.84   # HOTEL
    5  = 0.84   # OFFICE
    6  = 0.84   # RETAIL
    7  = 0.84   # FOODSTORE
    8  = 0.84   # RESTAURANT
    9  = 0.7    # INDUSTRIAL
    10 = 0.67   # SCHOOL
    11 = 0.84   # HOSPITAL
    12 = 0.67   # GYM
    13 = 0      # SWIMMING
    14 = 1      # SERVERROOM
    15 = 0      # PARKING
    16 = 1      # COOLROOM
    17 = 0.67   # LAB
    18 = 0.67   # MUSEUM
    19 = 0.67   # LIBRARY
}

foreach ($row in $hsValues.Keys) {
    $ws.Range("E$row").Value = $hsValues[$row]
}

$ws.Range("E12").Select()
